$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 9.258916175188528
$ws.Range("C2").Value = 5.440462225162122
$ws.Range("E2").Value = 16.4334201268971
$ws.Range("F2").Value = 16.86991607391245
$ws.Range("G2").Value = 30.99461872526152
$ws.Range("H2").Value = 14.74971678999694
$ws.Range("I2").Value = 21.17683375297929
$ws.Range("K2").Value = 9.236828046466716

$ws.Range("B3").Value = 8.92679331772303
$ws.Range("C3").Value = 5.134884828657177
$ws.Range("E3").Value = 15.50512914951802
$ws.Range("F3").Value = 15.89584955866815
$ws.Range("G3").Value = 31.0175067838662
$ws.Range("H3").Value = 14.80764593051684
$ws.Range("I3").Value = 21.28059472904867
$ws.Range("K3").Value = 9.004434903847605

$ws.Range("B4").Value = 8.718203778278932
$ws.Range("C4").Value = 4.936747720004851
$ws.Range("E4").Value = 14.91084546013573
$ws.Range("F4").Value = 15.26997757108489
$ws.Range("G4").Value = 31.04507282986189
$ws.Range("H4").Value = 14.84641504118384
$ws.Range("I4").Value = 21.34961050321206
$ws.Range("K4").Value = 8.860603488149579

$ws.Range("B5").Value = 8.632167013351271
$ws.Range("C5").Value = 4.853387736614263
$ws.Range("E5").Value = 14.66282873928581
$ws.Range("F5").Value = 15.008197319934
$ws.Range("G5").Value = 31.059679292841
$ws.Range("H5").Value = 14.86301583165411
$ws.Range("I5").Value = 21.37906361064209
$ws.Range("K5").Value = 8.801794795587121

$ws.Range("B6").Value = 8.617822463206073
$ws.Range("C6").Value = 4.839388666147411
$ws.Range("E6").Value = 14.62130193024104
$ws.Range("F6").Value = 14.96433081551589
$ws.Range("G6").Value = 31.06230762341697
$ws.Range("H6").Value = 14.86582074399073
$ws.Range("I6").Value = 21.3840343579081
$ws.Range("K6").Value = 8.792020532492135

$ws.Range("B7").Value = 8.717047454355546
$ws.Range("C7").Value = 4.935634057496239
$ws.Range("E7").Value = 14.90752387342068
$ws.Range("F7").Value = 15.26647399323133
$ws.Range("G7").Value = 31.0452561955862
$ws.Range("H7").Value = 14.8466356813557
$ws.Range("I7").Value = 21.3500023461065
$ws.Range("K7").Value = 8.859811045354915

$ws.Range("B8").Value = 9.145447283299962
$ws.Range("C8").Value = 5.337294188863268
$ws.Range("E8").Value = 16.11853800994067
$ws.Range("F8").Value = 16.5399640634477
$ws.Range("G8").Value = 30.99969203256501
$ws.Range("H8").Value = 14.76902489560895
$ws.Range("I8").Value = 21.21150584602047
$ws.Range("K8").Value = 9.156988893202399

$ws.Range("B9").Value = 9.943107946325384
$ws.Range("C9").Value = 6.040724438320705
$ws.Range("E9").Value = 18.35640402472106
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 31.01850795970317
$ws.Range("H9").Value = 14.64234479719333
$ws.Range("I9").Value = 20.98227052661071
$ws.Range("K9").Value = 9.726930725718466

$ws.Range("B10").Value = 10.49691408335786
$ws.Range("C10").Value = 6.505436048243415
$ws.Range("E10").Value = 19.99034357474723
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 31.09933297816438
$ws.Range("H10").Value = 14.5649927336569
$ws.Range("I10").Value = 20.84002347913778
$ws.Range("K10").Value = 10.13314365927513

$ws.Range("B11").Value = 10.74074820632179
$ws.Range("C11").Value = 6.70544981628128
$ws.Range("E11").Value = 20.69180957737874
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 31.15079974870119
$ws.Range("H11").Value = 14.53325227682818
$ws.Range("I11").Value = 20.78107120348628
$ws.Range("K11").Value = 10.31433099398112

$ws.Range("B12").Value = 10.83184017971648
$ws.Range("C12").Value = 6.779548876126569
$ws.Range("E12").Value = 20.95145958777231
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 31.17240863443853
$ws.Range("H12").Value = 14.52173168712813
$ws.Range("I12").Value = 20.75958178901806
$ws.Range("K12").Value = 10.38235741705107

$ws.Range("B13").Value = 10.81227831099759
$ws.Range("C13").Value = 6.763663386685812
$ws.Range("E13").Value = 20.89580473411704
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 31.16766045693994
$ws.Range("H13").Value = 14.52419061718883
$ws.Range("I13").Value = 20.76417269496173
$ws.Range("K13").Value = 10.36773375561107

$ws.Range("B14").Value = 10.74826771903131
$ws.Range("C14").Value = 6.711578913023986
$ws.Range("E14").Value = 20.71329079345365
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 31.15253502687822
$ws.Range("H14").Value = 14.5322944558883
$ws.Range("I14").Value = 20.7792864923623
$ws.Range("K14").Value = 10.31993966545288

$ws.Range("B15").Value = 10.70889537993276
$ws.Range("C15").Value = 6.679461759896875
$ws.Range("E15").Value = 20.60071797270846
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 31.14354639530995
$ws.Range("H15").Value = 14.53732334576613
$ws.Range("I15").Value = 20.78865301137028
$ws.Range("K15").Value = 10.29058627611805

$ws.Range("B16").Value = 10.48080910347944
$ws.Range("C16").Value = 6.49213512503735
$ws.Range("E16").Value = 19.94366247413125
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 31.09626587577805
$ws.Range("H16").Value = 14.5671366827049
$ws.Range("I16").Value = 20.84399256163093
$ws.Range("K16").Value = 10.12122428477579

$ws.Range("B17").Value = 10.33875303900493
$ws.Range("C17").Value = 6.374296136719682
$ws.Range("E17").Value = 19.52989296775868
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 31.07103100399868
$ws.Range("H17").Value = 14.58631136658815
$ws.Range("I17").Value = 20.87942076848974
$ws.Range("K17").Value = 10.01635411209758

$ws.Range("B18").Value = 10.25628823512803
$ws.Range("C18").Value = 6.305446684408095
$ws.Range("E18").Value = 19.28796526047643
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 31.05790002171323
$ws.Range("H18").Value = 14.59766443629974
$ws.Range("I18").Value = 20.90033956952372
$ws.Range("K18").Value = 9.955700300529628

$ws.Range("B19").Value = 10.22823945917488
$ws.Range("C19").Value = 6.281951573523608
$ws.Range("E19").Value = 19.20537522934778
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 31.05369146261555
$ws.Range("H19").Value = 14.60156399468826
$ws.Range("I19").Value = 20.90751510353936
$ws.Range("K19").Value = 9.935108528092123

$ws.Range("B20").Value = 10.35395419218923
$ws.Range("C20").Value = 6.386951243059115
$ws.Range("E20").Value = 19.57434680472517
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 31.07357406239926
$ws.Range("H20").Value = 14.58423660250364
$ws.Range("I20").Value = 20.87559329033195
$ws.Range("K20").Value = 10.0275529258282

$ws.Range("B21").Value = 10.76710348626671
$ws.Range("C21").Value = 6.726921968992042
$ws.Range("E21").Value = 20.76706158673494
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 31.15692018205038
$ws.Range("H21").Value = 14.52990059915311
$ws.Range("I21").Value = 20.77482450074913
$ws.Range("K21").Value = 10.33399434290955

$ws.Range("B22").Value = 11.02984725386046
$ws.Range("C22").Value = 6.939542126837744
$ws.Range("E22").Value = 21.51172955615892
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 31.22374585163042
$ws.Range("H22").Value = 14.49729789179474
$ws.Range("I22").Value = 20.71383373125906
$ws.Range("K22").Value = 10.53083263379392

$ws.Range("B23").Value = 10.89030492160418
$ws.Range("C23").Value = 6.826939379452428
$ws.Range("E23").Value = 21.11746296831373
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 31.18694843475225
$ws.Range("H23").Value = 14.51443137223988
$ws.Range("I23").Value = 20.74593804018281
$ws.Range("K23").Value = 10.42611184643053

$ws.Range("B24").Value = 10.34708422061039
$ws.Range("C24").Value = 6.381233302500677
$ws.Range("E24").Value = 19.55426183114713
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 31.07242005705047
$ws.Range("H24").Value = 14.58517357775971
$ws.Range("I24").Value = 20.87732197831225
$ws.Range("K24").Value = 10.02249106618323

$ws.Range("B25").Value = 9.732566979830763
$ws.Range("C25").Value = 5.859538225752008
$ws.Range("E25").Value = 17.72607729076346
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 31.00171396805856
$ws.Range("H25").Value = 14.67386653903453
$ws.Range("I25").Value = 21.03971427193897
$ws.Range("K25").Value = 9.57462517921706
